$d = $word.ActiveDocument

$replacements = @(
    @("41×67=2747", "34×27=918"),
    @("38×62=2356", "85×54=4590"),
    @("71×70=4970", "35×72=2520"),
    @("81×83=6723", "12×11=132"),
    @("71×34=2414", "37×33=1221"),
    @("16×89=1424", "27×17=459"),
    @("51×73=3723", "31×73=2263"),
    @("23×65=1495", "31×30=930"),
    @("94×31=2914", "88×11=968"),
    @("42×53=2226", "89×35=3115"),
    @("37×76=2812", "69×96=6624"),
    @("74×28=2072", "90×42=3780"),
    @("23×70=1610", "12×40=480"),
    @("66×95=6270", "95×32=3040"),
    @("89×95=8455", "65×28=1820"),
    @("61×63=3843", "14×25=350"),
    @("44×44=1936", "63×92=5796"),
    @("64×67=4288", "57×27=1539"),
    @("81×68=5508", "81×40=3240"),
    @("68×49=3332", "97×51=4947"),
    @("18×89=1602", "16×33=528"),
    @("27×82=2214", "47×59=2773"),
    @("11×91=1001", "37×19=703"),
    @("69×44=3036", "74×91=6734"),
    @("71×33=2343", "13×75=975"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
